$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two measured values in row 2 (B2: Min, C2: Max)
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 11

# Move/save the active cell selection to C4 (matches the saved view state)
$ws.Activate()
$ws.Range("C4").Select()
